# Weekly update: insert a new daily price record for "Ají" (Feria Lagunitas
# de Puerto Montt) as row 126, pushing the existing history (old rows
# 126-148) down by one row to 127-149.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 126; everything from 126-148 shifts to 127-149
# and the sheet's used-range / dimension grows to R149 automatically.
$ws.Rows(126).Insert()

# Populate the newly inserted row 126 with the new record.
$ws.Cells.Item(126, 1).Value  = 4
$ws.Cells.Item(126, 2).Value  = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(126, 3).Value  = "Los Lagos"
$ws.Cells.Item(126, 4).Value  = 44474
$ws.Cells.Item(126, 5).Value  = 10
$ws.Cells.Item(126, 6).Value  = 100112021
$ws.Cells.Item(126, 7).Value  = "Ají"
$ws.Cells.Item(126, 8).Value  = "Inferno"
$ws.Cells.Item(126, 9).Value  = "Primera"
$ws.Cells.Item(126, 10).Value = 140
$ws.Cells.Item(126, 11).Value = 47000
$ws.Cells.Item(126, 12).Value = 50000
$ws.Cells.Item(126, 13).Value = 48500
$ws.Cells.Item(126, 14).Value = "$/caja 12 kilos"
$ws.Cells.Item(126, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(126, 16).Value = 4042
$ws.Cells.Item(126, 17).Value = 12
$ws.Cells.Item(126, 18).Value = "Hortaliza"
